$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Your Clusters for Final Project" (column D) values for each
# student row — these were blank / placeholder ("0 or 1", "…") before and
# are now finalized 0/1 indicator values.
$dValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    14 = 1
    15 = 0
    16 = 1
    17 = 0
    18 = 0
}

foreach ($row in $dValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}

# Update the view: scroll so column B is the left-most visible column, and
# move the active selection to D12.
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D12").Select()
